$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.176.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +7.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.585.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +7.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9895"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "297.39"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3604"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3325"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.11"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.112"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06902"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.0000"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.25"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.782"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.494"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9902"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.587.33"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001057"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06556"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +10.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.69"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +9.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.73"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +8.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.883"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.49"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.185.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.361"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.486"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +16.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.07"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +11.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.757.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.25"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.899"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.790"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +17.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9100"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +13.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08107"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.626"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.62"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +11.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.070"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.31%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05962"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.263"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +11.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02164"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.43%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1966"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.50%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9903"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5732"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.747"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.96%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5539"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.26%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.22"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.925"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.95%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.03"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.08%  "
